$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 126, shifting the existing rows 126-129 down to 127-130.
$ws.Rows.Item(126).Insert()

# Populate the newly inserted row 126 with the new weekly price record
# (same market/category/variety/quality/origin as the row it was inserted above of,
# but with its own date, volume and price figures).
$ws.Cells.Item(126, 1).Value = 7
$ws.Cells.Item(126, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(126, 3).Value = "Ñuble"
$ws.Cells.Item(126, 4).Value = 45239
$ws.Cells.Item(126, 5).Value = 16
$ws.Cells.Item(126, 6).Value = 100112001
$ws.Cells.Item(126, 7).Value = "Berenjena"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 30
$ws.Cells.Item(126, 11).Value = 12000
$ws.Cells.Item(126, 12).Value = 12000
$ws.Cells.Item(126, 13).Value = 12000
$ws.Cells.Item(126, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(126, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(126, 16).Value = 200
$ws.Cells.Item(126, 17).Value = 60
$ws.Cells.Item(126, 18).Value = "Hortaliza"
